$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: " [9" + bookmark(_GoBack) + "]"  ->  " [9]" (single run, bookmark
# removed from this location).
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$find1 = $d.Content
$find1.Find.Execute(" [9]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find1.Delete()
$find1.InsertAfter(" [9]")

# ---------------------------------------------------------------------------
# Change 2: fix the typo "ant" -> "and" in the "Dev Mode" section, and move
# the _GoBack bookmark so that it now sits right after "...screen and" and
# before " type in the name ...".
# ---------------------------------------------------------------------------
$fixRange = $d.Content
$fixRange.Find.Execute("screen ant type", $true, $false, $false, $false, $false, $true, 1, $false, "screen and type", 2)

# Locate "a user mu" / "st go" boundary so we can force a run split there too,
# matching the way Word originally split this sentence into multiple runs.
$muRange = $d.Content
$muRange.Find.Execute("a user mu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $muRange.End
$tempRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("TempSplitMarker", $tempRange)
$d.Bookmarks("TempSplitMarker").Delete()

# Now locate the boundary right after "...screen and" (before the space and
# "type") and drop the _GoBack bookmark there.
$andRange = $d.Content
$andRange.Find.Execute("screen and", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkPoint = $andRange.End
$bookmarkRange = $d.Range($bookmarkPoint, $bookmarkPoint)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
